$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.942.21"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.10%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.510.97"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.66%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "570.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.07%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "185.05"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.90%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.503.92"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.75%  "

$ws.Range("E8").Value = "  -2.73%  "

$ws.Range("E9").Value = "  +0.07%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.186"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.39%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.651"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.30%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.20"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.47%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000301"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.54%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.44"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.77%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.073.01"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.77%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.35"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.11%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.508.09"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.81%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.907.43"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.11%  "

$ws.Range("E19").Value = "  -2.82%  "

$ws.Range("E20").Value = "  -1.13%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "540.55"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +14.12%  "

$ws.Range("E22").Value = "  -2.80%  "

$ws.Range("E23").Value = "  -4.68%  "

$ws.Range("E24").Value = "  -0.59%  "

$ws.Range("E25").Value = "  +0.74%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "93.97"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.16%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.92"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.48%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.80"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.85%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.11"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.32%  "

$ws.Range("E30").Value = "  -1.16%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.24"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.00%  "

$ws.Range("E32").Value = "  +3.05%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "64.81"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.01%  "

$ws.Range("E34").Value = "  -3.84%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "560.50"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.59%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "37.89"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.39%  "

$ws.Range("E37").Value = "  +0.18%  "

$ws.Range("E38").Value = "  +0.76%  "

$ws.Range("E39").Value = "  +4.34%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0764"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.86%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.15"
$ws.Range("D41").Style = "Normal"

$ws.Range("E42").Value = "  -3.07%  "

$ws.Range("E43").Value = "  -2.76%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.55"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.58%  "

$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.234.23"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.24%  "

$ws.Range("B46").Value = "ThetaToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.99"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.78%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0439"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.54%  "

$ws.Range("E48").Value = "  -2.47%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.93"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.03%  "

$ws.Range("E50").Value = "  -0.11%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "137.70"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.93%  "
